$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($r=36; $r -le 41; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    Write-Host "Row $r A=$a B=$b C=$c"
}
